# Update the "Metabolites" sheet with results from the new Ecoli model runs
# (100s and 300s simulations). Several metabolite concentrations that used
# to be plain 0 now hold the high-precision numeric results of the model,
# stored as text (as the original workbook already did for a couple of
# cells), while DPG goes back to a plain 0.
#
# NOTE: assigning a numeric-looking string directly to Range.Value makes
# Excel reinterpret it as a number (and, with a leading apostrophe, tags the
# cell with a "quote prefix" style instead of reusing the sheet's existing
# styles). Building the text through a formula and then pasting back only
# the values keeps the cell as a genuine shared-string text value without
# disturbing its existing number format / style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metabolites")

function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}

Set-TextValue "B18" "8.7363521885873"
Set-TextValue "B23" "2.7926086603330034"
Set-TextValue "B36" "0.00010000000001936753"
Set-TextValue "B43" "2.676318153277067"
Set-TextValue "B45" "0.05248431161643879"
Set-TextValue "B49" "14.008413729564946"
Set-TextValue "B53" "1.315108425425397"

# DPG no longer has a measured value; reset it back to a plain numeric 0.
$ws.Range("B25").Value = 0

$excel.CutCopyMode = $false
